$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 97 (pushing the existing rows 97-177
# down to 98-178, which is exactly what the canonical-XML diff shows: every
# row from 97 onward is replaced by the row that used to precede it, and a
# new final row 178 appears carrying the data that used to live in row 177).
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with its data.
$ws.Range("A97").Value = 11
$ws.Range("B97").Value = "Vega Monumental Concepción"
$ws.Range("C97").Value = "Bíobío"
$ws.Range("D97").Value = 44729
$ws.Range("E97").Value = 8
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100108
$ws.Range("H97").Value = "Tropicales y subtropicales"
$ws.Range("I97").Value = 100108005
$ws.Range("J97").Value = "Piña"
$ws.Range("K97").Value = "Caramelo"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 200
$ws.Range("N97").Value = 16000
$ws.Range("O97").Value = 17000
$ws.Range("P97").Value = 16500
$ws.Range("Q97").Value = "$/caja 12 unidades"
$ws.Range("R97").Value = "Ecuador"
$ws.Range("S97").Value = 1375
$ws.Range("T97").Value = 12
